$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.367.44"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  -2.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.189.55"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  -3.68%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.47"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  -1.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.85"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  -5.08%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.186.26"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  -3.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "  -3.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  -4.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.27"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  -4.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.451"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "  -4.82%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000236"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  -5.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.32"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  -4.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.708.00"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  -3.90%  "

$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.195.13"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  -3.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.415.73"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  -2.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.60"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  -4.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "458.76"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  -4.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.92"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  -2.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.707"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  -4.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.63"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  -4.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.35"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "  -1.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.92"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "  -2.12%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "  -3.36%  "

$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.86"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  -4.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.91"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  -5.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.05"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  -5.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.22"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  -6.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.102"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  -3.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.38"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  -6.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.04"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  -5.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.82"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  -2.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.04"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  -4.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0691"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  -8.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0386"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = "  -3.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.972.94"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  -2.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "411.06"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  -4.77%  "

$ws.Range("B42").Value = "Kaspa"

$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  +2.33%  "

$ws.Range("B43").Value = "dogwifhat"

$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.66"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  -4.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.03"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  -4.91%  "

$ws.Range("B45").Value = "TheGraph"

$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.250"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  -6.84%  "

$ws.Range("B46").Value = "Fetch.AI"

$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  -3.37%  "

$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.83"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  +1.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.66"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  -3.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.94"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  +0.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.111"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = "  -3.94%  "
